$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B1").Value = "Digitale_Woche"
$ws.Range("B1").Select()
